# Auto-generated: apply updated profit-calculation values to Sheets per scheduled-runner diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 38.166668
$ws.Range("I4").Value = 38.166668
$ws.Range("K4").Value = 38.166668
$ws.Range("M4").Value = 75.833332
# Row 40
$ws.Range("H40").Value = 21405.234
$ws.Range("J40").Value = 31497.5
$ws.Range("L40").Value = 31497.5
$ws.Range("N40").Value = -31847.5
# Row 41
$ws.Range("H41").Value = 4724.75
$ws.Range("I41").Value = 3500
$ws.Range("J41").Value = 5133
$ws.Range("K41").Value = 3500
$ws.Range("L41").Value = 5133
$ws.Range("N41").Value = -6013
$ws.Range("M41").Value = -3060
# Row 70
$ws.Range("H70").Value = 5286.5454
$ws.Range("J70").Value = 5944.4443
$ws.Range("L70").Value = 17833.3329
$ws.Range("N70").Value = -18373.3329
# Row 73
$ws.Range("H73").Value = 5286.5454
$ws.Range("J73").Value = 5944.4443
$ws.Range("L73").Value = 17833.3329
$ws.Range("N73").Value = -19705.3329
# Row 92
$ws.Range("H92").Value = 592.5333000000001
$ws.Range("I92").Value = 498.7857
$ws.Range("J92").Value = 1905
$ws.Range("K92").Value = 498.7857
$ws.Range("L92").Value = 1905
$ws.Range("M92").Value = 749.2143
$ws.Range("N92").Value = -4401
# Row 96
$ws.Range("H96").Value = 1361.75
$ws.Range("J96").Value = 2303.875
$ws.Range("L96").Value = 6911.625
$ws.Range("N96").Value = -9657.625
# Row 98
$ws.Range("H98").Value = 1475.2222
$ws.Range("I98").Value = 1370.3334
$ws.Range("K98").Value = 1370.3334
$ws.Range("M98").Value = 127.6666
# Row 122
$ws.Range("H122").Value = 1475.2222
$ws.Range("I122").Value = 1370.3334
$ws.Range("K122").Value = 4111.0002
$ws.Range("M122").Value = -1661.0002
# Row 132
$ws.Range("H132").Value = 12669.384
$ws.Range("I132").Value = 2282.1428
$ws.Range("K132").Value = 6846.428400000001
$ws.Range("M132").Value = -4316.428400000001
# Row 135
$ws.Range("H135").Value = 2728.8635
$ws.Range("J135").Value = 7007
$ws.Range("L135").Value = 63063
$ws.Range("N135").Value = -68133
# Row 137
$ws.Range("H137").Value = 10421902
$ws.Range("I137").Value = 2496.4
$ws.Range("K137").Value = 7489.200000000001
$ws.Range("M137").Value = -4939.200000000001
# Row 141
$ws.Range("H141").Value = 11323.5
$ws.Range("I141").Value = 12441.143
$ws.Range("K141").Value = 37323.429
$ws.Range("M141").Value = -32143.429

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3592.218
$ws.Range("I32").Value = 2019.2924
$ws.Range("J32").Value = 11456.846
$ws.Range("K32").Value = 2019.2924
$ws.Range("L32").Value = 11456.846
$ws.Range("M32").Value = -1732.2924
$ws.Range("N32").Value = -12030.846
# Row 33
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
# Row 61
$ws.Range("H61").Value = 12589.294
$ws.Range("I61").Value = 15540.77
$ws.Range("K61").Value = 15540.77
$ws.Range("M61").Value = -15328.77
# Row 74
$ws.Range("H74").Value = 17859620
$ws.Range("I74").Value = 62501500
$ws.Range("K74").Value = 62501500
$ws.Range("M74").Value = -62500626
# Row 77
$ws.Range("H77").Value = 17859620
$ws.Range("I77").Value = 62501500
$ws.Range("K77").Value = 312507500
$ws.Range("M77").Value = -312503132
# Row 110
$ws.Range("H110").Value = 4418.2144
$ws.Range("I110").Value = 2851.4443
$ws.Range("K110").Value = 2851.4443
$ws.Range("M110").Value = -806.4443000000001
# Row 136
$ws.Range("H136").Value = 12589.294
$ws.Range("I136").Value = 15540.77
$ws.Range("K136").Value = 46622.31
$ws.Range("M136").Value = -44072.31

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 869124.9
$ws.Range("J99").Value = 1325.8334
$ws.Range("L99").Value = 1325.8334
$ws.Range("N99").Value = -4321.8334
# Row 134
$ws.Range("H134").Value = 7426.2856
$ws.Range("I134").Value = 5333.3335
$ws.Range("J134").Value = 8996
$ws.Range("K134").Value = 16000.0005
$ws.Range("L134").Value = 26988
$ws.Range("M134").Value = -13465.0005
$ws.Range("N134").Value = -32058

$ws = $wb.Worksheets.Item("CRP")
# Row 33
$ws.Range("H33").Value = 8999
$ws.Range("I33").Value = 8999
$ws.Range("K33").Value = 8999
$ws.Range("M33").Value = -8620
# Row 132
$ws.Range("H132").Value = 50022884
$ws.Range("I132").Value = 62515540
$ws.Range("K132").Value = 187546620
$ws.Range("M132").Value = -187544090
# Row 134
$ws.Range("H134").Value = 2484.2632
$ws.Range("I134").Value = 2472.2778
$ws.Range("K134").Value = 7416.8334
$ws.Range("M134").Value = -4881.8334
# Row 140
$ws.Range("H140").Value = 94910
$ws.Range("J140").Value = 94910
$ws.Range("L140").Value = 94910
$ws.Range("N140").Value = -105270

$ws = $wb.Worksheets.Item("CUL")
# Row 96
$ws.Range("H96").Value = 12998.833
$ws.Range("I96").Value = 12997
$ws.Range("K96").Value = 38991
$ws.Range("M96").Value = -36932
# Row 129
$ws.Range("H129").Value = 3500
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 267.8
$ws.Range("I97").Value = 238.5
$ws.Range("J97").Value = 326.4
$ws.Range("K97").Value = 238.5
$ws.Range("L97").Value = 326.4
$ws.Range("M97").Value = 257.5
$ws.Range("N97").Value = -1318.4
# Row 132
$ws.Range("H132").Value = 6609.1816
$ws.Range("I132").Value = 7877.6
$ws.Range("J132").Value = 5552.1665
$ws.Range("K132").Value = 23632.8
$ws.Range("L132").Value = 16656.4995
$ws.Range("M132").Value = -21102.8
$ws.Range("N132").Value = -21716.4995

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1776.3334
$ws.Range("I16").Value = 1452.3572
$ws.Range("J16").Value = 2424.2856
$ws.Range("K16").Value = 1452.3572
$ws.Range("L16").Value = 2424.2856
$ws.Range("M16").Value = -1282.3572
$ws.Range("N16").Value = -2764.2856
# Row 22
$ws.Range("H22").Value = 1178.3572
$ws.Range("I22").Value = 999.875
$ws.Range("J22").Value = 1416.3334
$ws.Range("K22").Value = 999.875
$ws.Range("L22").Value = 1416.3334
$ws.Range("M22").Value = -704.875
$ws.Range("N22").Value = -2006.3334
# Row 27
$ws.Range("H27").Value = 1178.3572
$ws.Range("I27").Value = 999.875
$ws.Range("J27").Value = 1416.3334
$ws.Range("K27").Value = 999.875
$ws.Range("L27").Value = 1416.3334
$ws.Range("M27").Value = -892.875
$ws.Range("N27").Value = -1630.3334
# Row 46
$ws.Range("H46").Value = 6409.2856
$ws.Range("I46").Value = 2037
$ws.Range("J46").Value = 7138
$ws.Range("K46").Value = 2037
$ws.Range("L46").Value = 7138
$ws.Range("M46").Value = -1849
$ws.Range("N46").Value = -7514
# Row 55
$ws.Range("H55").Value = 510.42856
$ws.Range("I55").Value = 89
$ws.Range("K55").Value = 89
$ws.Range("M55").Value = 84
# Row 82
$ws.Range("H82").Value = 2232676.2
$ws.Range("I82").Value = 2841297
$ws.Range("J82").Value = 1066.3334
$ws.Range("K82").Value = 2841297
$ws.Range("L82").Value = 1066.3334
$ws.Range("M82").Value = -2840936
$ws.Range("N82").Value = -1788.3334
# Row 85
$ws.Range("H85").Value = 2232676.2
$ws.Range("I85").Value = 2841297
$ws.Range("J85").Value = 1066.3334
$ws.Range("K85").Value = 2841297
$ws.Range("L85").Value = 1066.3334
$ws.Range("M85").Value = -2840049
$ws.Range("N85").Value = -3562.3334
# Row 132
$ws.Range("H132").Value = 4638.9565
$ws.Range("I132").Value = 4278
$ws.Range("J132").Value = 5315.75
$ws.Range("K132").Value = 12834
$ws.Range("L132").Value = 15947.25
$ws.Range("M132").Value = -10304
$ws.Range("N132").Value = -21007.25

$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 10307.085
$ws.Range("I136").Value = 5889.294
$ws.Range("J136").Value = 12095.238
$ws.Range("K136").Value = 17667.882
$ws.Range("L136").Value = 36285.714
$ws.Range("M136").Value = -15117.882
$ws.Range("N136").Value = -41385.714
